$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the stale "ParticipantsTab" Neo4j query (B2) with the updated
#     query text (adds diagnosis/genomic_info lookups and sorts the sample
#     list via apoc.coll.sort). Setting the cell value causes the now-unused
#     old shared string to be dropped and the new text appended, matching
#     the authoring diff's sharedStrings churn automatically.
$newParticipantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['PDF']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$ws.Cells.Item(2, 2).Value = $newParticipantQuery

# --- The longer replacement text wraps across more lines, so row 2 grows
#     taller to keep the whole query visible.
$ws.Rows.Item(2).RowHeight = 279

# --- Scroll the view down a bit and move the active selection to B4, as in
#     the saved view state of the edited workbook.
$ws.Range("B4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
